$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (601) down across the new rows (602:613)
$ws.Range("A601:V601").Copy()
$ws.Range("A602:V613").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 602
$ws.Range("A602").Value = "Entrainement"
$ws.Range("B602").Value = 45932
$ws.Range("C602").Value = "Global"
$ws.Range("D602").Value = "J-2"
$ws.Range("E602").Value = "Omar Benyounes"
$ws.Range("F602").Value = "center midfield"
$ws.Range("G602").Value = "01:20:28"
$ws.Range("H602").Value = 5.05
$ws.Range("I602").Value = 0.18
$ws.Range("J602").Value = 4.86
$ws.Range("K602").Value = 0.18
$ws.Range("L602").Value = 0.01
$ws.Range("M602").Value = 0
$ws.Range("N602").Value = 0
$ws.Range("O602").Value = 0
$ws.Range("P602").Value = 3.68
$ws.Range("Q602").Value = 21.93
$ws.Range("R602").Value = 4.59
$ws.Range("S602").Value = 16
$ws.Range("T602").Value = 7
$ws.Range("U602").Value = 17
$ws.Range("V602").Value = 1

# Row 603
$ws.Range("A603").Value = "Entrainement"
$ws.Range("B603").Value = 45932
$ws.Range("C603").Value = "Global"
$ws.Range("D603").Value = "J-2"
$ws.Range("E603").Value = "Mattheo Haon"
$ws.Range("F603").Value = "right back"
$ws.Range("G603").Value = "01:21:01"
$ws.Range("H603").Value = 5.4
$ws.Range("I603").Value = 0.19
$ws.Range("J603").Value = 5.21
$ws.Range("K603").Value = 0.13
$ws.Range("L603").Value = 0.05
$ws.Range("M603").Value = 0.01
$ws.Range("N603").Value = 0
$ws.Range("O603").Value = 1
$ws.Range("P603").Value = 3.89
$ws.Range("Q603").Value = 26.37
$ws.Range("R603").Value = 4.48
$ws.Range("S603").Value = 23
$ws.Range("T603").Value = 2
$ws.Range("U603").Value = 9
$ws.Range("V603").Value = 6

# Row 604
$ws.Range("A604").Value = "Entrainement"
$ws.Range("B604").Value = 45932
$ws.Range("C604").Value = "Global"
$ws.Range("D604").Value = "J-2"
$ws.Range("E604").Value = "Hedi Nasri"
$ws.Range("F604").Value = "right back"
$ws.Range("G604").Value = "01:10:56"
$ws.Range("H604").Value = 4.4
$ws.Range("I604").Value = 0.11
$ws.Range("J604").Value = 4.28
$ws.Range("K604").Value = 0.07
$ws.Range("L604").Value = 0.02
$ws.Range("M604").Value = 0.02
$ws.Range("N604").Value = 0
$ws.Range("O604").Value = 2
$ws.Range("P604").Value = 2.88
$ws.Range("Q604").Value = 28.76
$ws.Range("R604").Value = 4.8
$ws.Range("S604").Value = 17
$ws.Range("T604").Value = 10
$ws.Range("U604").Value = 14
$ws.Range("V604").Value = 6

# Row 605
$ws.Range("A605").Value = "Entrainement"
$ws.Range("B605").Value = 45933
$ws.Range("C605").Value = "Global"
$ws.Range("D605").Value = "J-1"
$ws.Range("E605").Value = "Mattheo Haon"
$ws.Range("F605").Value = "right back"
$ws.Range("G605").Value = "01:17:41"
$ws.Range("H605").Value = 6.03
$ws.Range("I605").Value = 0.45
$ws.Range("J605").Value = 5.56
$ws.Range("K605").Value = 0.3
$ws.Range("L605").Value = 0.15
$ws.Range("M605").Value = 0.02
$ws.Range("N605").Value = 0
$ws.Range("O605").Value = 3
$ws.Range("P605").Value = 4.59
$ws.Range("Q605").Value = 28.97
$ws.Range("R605").Value = 5.1
$ws.Range("S605").Value = 54
$ws.Range("T605").Value = 13
$ws.Range("U605").Value = 34
$ws.Range("V605").Value = 6

# Row 606
$ws.Range("A606").Value = "Entrainement"
$ws.Range("B606").Value = 45933
$ws.Range("C606").Value = "Global"
$ws.Range("D606").Value = "J-1"
$ws.Range("E606").Value = "Sofiane Belle"
$ws.Range("F606").Value = "left forward"
$ws.Range("G606").Value = "01:20:46"
$ws.Range("H606").Value = 5.24
$ws.Range("I606").Value = 0.28
$ws.Range("J606").Value = 4.94
$ws.Range("K606").Value = 0.24
$ws.Range("L606").Value = 0.05
$ws.Range("M606").Value = 0
$ws.Range("N606").Value = 0
$ws.Range("O606").Value = 1
$ws.Range("P606").Value = 3.76
$ws.Range("Q606").Value = 25.04
$ws.Range("R606").Value = 4.32
$ws.Range("S606").Value = 16
$ws.Range("T606").Value = 2
$ws.Range("U606").Value = 19
$ws.Range("V606").Value = 6

# Row 607
$ws.Range("A607").Value = "Entrainement"
$ws.Range("B607").Value = 45933
$ws.Range("C607").Value = "Global"
$ws.Range("D607").Value = "J-1"
$ws.Range("E607").Value = "Kamal Bafounta"
$ws.Range("F607").Value = "center midfield"
$ws.Range("G607").Value = "01:19:03"
$ws.Range("H607").Value = 5.93
$ws.Range("I607").Value = 0.28
$ws.Range("J607").Value = 5.64
$ws.Range("K607").Value = 0.21
$ws.Range("L607").Value = 0.08
$ws.Range("M607").Value = 0
$ws.Range("N607").Value = 0
$ws.Range("O607").Value = 0
$ws.Range("P607").Value = 4.46
$ws.Range("Q607").Value = 23.98
$ws.Range("R607").Value = 4.59
$ws.Range("S607").Value = 29
$ws.Range("T607").Value = 3
$ws.Range("U607").Value = 16
$ws.Range("V607").Value = 2

# Row 608
$ws.Range("A608").Value = "Entrainement"
$ws.Range("B608").Value = 45933
$ws.Range("C608").Value = "Global"
$ws.Range("D608").Value = "J-1"
$ws.Range("E608").Value = "Emmanuel Valey"
$ws.Range("F608").Value = "left forward"
$ws.Range("G608").Value = "01:19:49"
$ws.Range("H608").Value = 6.19
$ws.Range("I608").Value = 0.38
$ws.Range("J608").Value = 5.8
$ws.Range("K608").Value = 0.22
$ws.Range("L608").Value = 0.12
$ws.Range("M608").Value = 0.04
$ws.Range("N608").Value = 0
$ws.Range("O608").Value = 7
$ws.Range("P608").Value = 4.23
$ws.Range("Q608").Value = 27.93
$ws.Range("R608").Value = 5.24
$ws.Range("S608").Value = 46
$ws.Range("T608").Value = 12
$ws.Range("U608").Value = 40
$ws.Range("V608").Value = 12

# Row 609
$ws.Range("A609").Value = "Entrainement"
$ws.Range("B609").Value = 45933
$ws.Range("C609").Value = "Global"
$ws.Range("D609").Value = "J-1"
$ws.Range("E609").Value = "Omar Benyounes"
$ws.Range("F609").Value = "center midfield"
$ws.Range("G609").Value = "01:17:52"
$ws.Range("H609").Value = 5.61
$ws.Range("I609").Value = 0.37
$ws.Range("J609").Value = 5.23
$ws.Range("K609").Value = 0.25
$ws.Range("L609").Value = 0.13
$ws.Range("M609").Value = 0
$ws.Range("N609").Value = 0
$ws.Range("O609").Value = 1
$ws.Range("P609").Value = 4.22
$ws.Range("Q609").Value = 25.19
$ws.Range("R609").Value = 4.49
$ws.Range("S609").Value = 25
$ws.Range("T609").Value = 10
$ws.Range("U609").Value = 15
$ws.Range("V609").Value = 7

# Row 610
$ws.Range("A610").Value = "Entrainement"
$ws.Range("B610").Value = 45933
$ws.Range("C610").Value = "Global"
$ws.Range("D610").Value = "J-1"
$ws.Range("E610").Value = "Jeremie Laurent"
$ws.Range("F610").Value = "left forward"
$ws.Range("G610").Value = "01:20:39"
$ws.Range("H610").Value = 6.36
$ws.Range("I610").Value = 0.42
$ws.Range("J610").Value = 5.93
$ws.Range("K610").Value = 0.32
$ws.Range("L610").Value = 0.09
$ws.Range("M610").Value = 0.01
$ws.Range("N610").Value = 0
$ws.Range("O610").Value = 1
$ws.Range("P610").Value = 4.61
$ws.Range("Q610").Value = 27.36
$ws.Range("R610").Value = 5.25
$ws.Range("S610").Value = 52
$ws.Range("T610").Value = 10
$ws.Range("U610").Value = 34
$ws.Range("V610").Value = 8

# Row 611
$ws.Range("A611").Value = "Entrainement"
$ws.Range("B611").Value = 45933
$ws.Range("C611").Value = "Global"
$ws.Range("D611").Value = "J-1"
$ws.Range("E611").Value = "Karahali Souaré"
$ws.Range("F611").Value = "right forward"
$ws.Range("G611").Value = "01:18:45"
$ws.Range("H611").Value = 5.65
$ws.Range("I611").Value = 0.22
$ws.Range("J611").Value = 5.42
$ws.Range("K611").Value = 0.2
$ws.Range("L611").Value = 0.03
$ws.Range("M611").Value = 0
$ws.Range("N611").Value = 0
$ws.Range("O611").Value = 0
$ws.Range("P611").Value = 3.85
$ws.Range("Q611").Value = 24.47
$ws.Range("R611").Value = 5.49
$ws.Range("S611").Value = 64
$ws.Range("T611").Value = 18
$ws.Range("U611").Value = 33
$ws.Range("V611").Value = 15

# Row 612
$ws.Range("A612").Value = "Entrainement"
$ws.Range("B612").Value = 45933
$ws.Range("C612").Value = "Global"
$ws.Range("D612").Value = "J-1"
$ws.Range("E612").Value = "Levy Ndoutoume"
$ws.Range("F612").Value = "left back"
$ws.Range("G612").Value = "01:17:33"
$ws.Range("H612").Value = 5.31
$ws.Range("I612").Value = 0.29
$ws.Range("J612").Value = 5.01
$ws.Range("K612").Value = 0.19
$ws.Range("L612").Value = 0.11
$ws.Range("M612").Value = 0
$ws.Range("N612").Value = 0
$ws.Range("O612").Value = 1
$ws.Range("P612").Value = 3.97
$ws.Range("Q612").Value = 25.26
$ws.Range("R612").Value = 4.86
$ws.Range("S612").Value = 25
$ws.Range("T612").Value = 9
$ws.Range("U612").Value = 24
$ws.Range("V612").Value = 4

# Row 613
$ws.Range("A613").Value = "Entrainement"
$ws.Range("B613").Value = 45933
$ws.Range("C613").Value = "Global"
$ws.Range("D613").Value = "J-1"
$ws.Range("E613").Value = "Malik Boussaid"
$ws.Range("F613").Value = "right back"
$ws.Range("G613").Value = "01:18:26"
$ws.Range("H613").Value = 6.44
$ws.Range("I613").Value = 0.44
$ws.Range("J613").Value = 5.98
$ws.Range("K613").Value = 0.29
$ws.Range("L613").Value = 0.13
$ws.Range("M613").Value = 0.03
$ws.Range("N613").Value = 0
$ws.Range("O613").Value = 4
$ws.Range("P613").Value = 4.57
$ws.Range("Q613").Value = 28.58
$ws.Range("R613").Value = 4.99
$ws.Range("S613").Value = 60
$ws.Range("T613").Value = 15
$ws.Range("U613").Value = 48
$ws.Range("V613").Value = 10

# Update the active selection to match the latest edit position
$ws.Range("E619").Select()
